# ---------------------------------------------------------------------------
# Commit: "add temp one match, two match / prepare for two match compute"
#
# 1. Populate a small "two match" reference table on Sheet3 (w / d / l,
#    wdl_w / wdl_d / wdl_l, the spread_* formulas, point_w/point_d/point_l
#    and the ww,dw,lw groupings), sized with generous column widths and a
#    tall, word-wrapped header row.
# 2. Make Sheet3 the active sheet/tab (it was "all" before).
# 3. Leave a fresh selection behind on the "all" sheet where the user had
#    scrolled to while double-checking the one-match data.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsAll   = $wb.Worksheets.Item("all")
$ws3     = $wb.Worksheets.Item("Sheet3")

# --- Sheet3: column widths -------------------------------------------------
$ws3.Columns.Item(1).ColumnWidth = 54.5
$ws3.Columns.Item(2).ColumnWidth = 56.83333333333
$ws3.Columns.Item(3).ColumnWidth = 41.66666666667

# --- Sheet3: header row (w / d / l) ----------------------------------------
$ws3.Range("A1").Value = "w"
$ws3.Range("B1").Value = "d"
$ws3.Range("C1").Value = "l"
$ws3.Rows.Item(1).RowHeight = 15

# --- Sheet3: wdl totals row --------------------------------------------------
$ws3.Range("A2").Value = "wdl_w"
$ws3.Range("B2").Value = "wdl_d"
$ws3.Range("C2").Value = "wdl_l"

# --- Sheet3: point_w / point_l row (point_d filled back in afterwards) -----
$ws3.Range("A4").Value = "point_w……"
$ws3.Range("C4").Value = "point_l……"

# --- Sheet3: the ww,dw,lw / wd,dd,ld / wl,dl,ll grouping row ---------------
$ws3.Range("A5").Value = "ww,dw,lw"
$ws3.Range("B5").Value = "wd,dd,ld"
$ws3.Range("C5").Value = "wl,dl,ll"

# --- Sheet3: back-fill point_d ----------------------------------------------
$ws3.Range("B4").Value = "point_d……"

# --- Sheet3: spread_w / spread_d / spread_l formulas row -------------------
$ws3.Range("B3").Value = "spread_d=(point_w_1_0+point _w_2_1+`n                      point_w_3_2+point_w_other)"
$ws3.Range("A3").Value = "spread_w =wdl-`n                      (point_w_1_0,point _w_2_1,`n                       point_w_3_2,point_w_other)"
$ws3.Range("C3").Value = "spread_l=wdl_l +(point_d_0_0+  point_d_1_1….)`nor`nspread_l=wdl_l+wdl_d"

$ws3.Range("A3:C3").WrapText = $true
$ws3.Rows.Item(3).RowHeight = 60

# --- "all": leave the selection where the user scrolled to -----------------
$wsAll.Range("C42").Select()

# --- make Sheet3 the active tab, with A3 selected ---------------------------
$ws3.Activate()
$ws3.Range("A3").Select()

Write-Output "edit applied"
